# Update countries & provincias Spain
# This script:
#  1) Updates COVID numbers for several countries (general data refresh).
#  2) Because Nigeria's, Uruguay's and Islas Virgenes Britanicas' figures were
#     updated, those countries now outrank the neighbor they were previously
#     listed under (Ghana, Mali, Butan respectively), so the pair of rows is
#     swapped (country name + data move as a unit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 1366074
$ws.Range("C4").Value = 18765
$ws.Range("D4").Value = 255654
$ws.Range("E4").Value = 1029683
$ws.Range("G4").Value = 700
$ws.Range("H4").Value = 80737

# --- Row 10: Alemania ---
$ws.Range("B10").Value = 171879
$ws.Range("C10").Value = 555
$ws.Range("E10").Value = 19910
$ws.Range("G10").Value = 20
$ws.Range("H10").Value = 7569

# --- Row 11: Brasil ---
$ws.Range("D11").Value = 64957
$ws.Range("E11").Value = 86619

# --- Row 15: Canada ---
$ws.Range("D15").Value = 32096
$ws.Range("E15").Value = 31882

# --- Row 51: Chequia ---
$ws.Range("B51").Value = 8123
$ws.Range("C51").Value = 28
$ws.Range("D51").Value = 4474
$ws.Range("E51").Value = 3369
$ws.Range("G51").Value = 4
$ws.Range("H51").Value = 280

# --- Row 56: Argentina ---
$ws.Range("B56").Value = 6034
$ws.Range("C56").Value = 258
$ws.Range("E56").Value = 3972

# --- Rows 63/64: Ghana & Nigeria swap order (Nigeria's numbers updated) ---
$ws.Range("A63").Value = "Nigeria"
$ws.Range("B63").Value = 4399
$ws.Range("C63").Value = 248
$ws.Range("D63").Value = 778
$ws.Range("E63").Value = 3478
$ws.Range("F63").Value = 4
$ws.Range("G63").Value = 15
$ws.Range("H63").Value = 143

$ws.Range("A64").Value = "Ghana"
$ws.Range("B64").Value = 4263
$ws.Range("C64").Value = 0
$ws.Range("D64").Value = 378
$ws.Range("E64").Value = 3863
$ws.Range("F64").Value = 5
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 22

# --- Rows 114/115: Mali & Uruguay swap order (Uruguay's numbers updated) ---
$ws.Range("A114").Value = "Uruguay"
$ws.Range("B114").Value = 707
$ws.Range("C114").Value = 5
$ws.Range("D114").Value = 517
$ws.Range("E114").Value = 171
$ws.Range("F114").Value = 8
$ws.Range("G114").Value = 1
$ws.Range("H114").Value = 19

$ws.Range("A115").Value = "Mali"
$ws.Range("B115").Value = 704
$ws.Range("C115").Value = 12
$ws.Range("D115").Value = 351
$ws.Range("E115").Value = 315
$ws.Range("F115").Value = 0
$ws.Range("G115").Value = 1
$ws.Range("H115").Value = 38

# --- Rows 212/213: Islas Virgenes Britanicas & Butan swap order (figures unchanged) ---
$ws.Range("A212").Value = "Butan"
$ws.Range("B212").Value = 7
$ws.Range("C212").Value = 0
$ws.Range("D212").Value = 5
$ws.Range("E212").Value = 2
$ws.Range("F212").Value = 0
$ws.Range("G212").Value = 0
$ws.Range("H212").Value = 0

$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("B213").Value = 7
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 4
$ws.Range("E213").Value = 2
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 1
